$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template cell (row 464, column A) carries the date style (s="2") that new rows should reuse.
$dateStyleSource = $ws.Cells.Item(464, 1)

$dateStyleSource.Copy($ws.Cells.Item(465, 1))
$ws.Cells.Item(465, 1).Value = 44539
$ws.Cells.Item(465, 2).Value = 4
$ws.Cells.Item(465, 3).Value = 29
$ws.Cells.Item(465, 4).Value = 837.424198671672

$dateStyleSource.Copy($ws.Cells.Item(466, 1))
$ws.Cells.Item(466, 1).Value = 44540
$ws.Cells.Item(466, 2).Value = 5
$ws.Cells.Item(466, 3).Value = 25
$ws.Cells.Item(466, 4).Value = 721.917412647993

$dateStyleSource.Copy($ws.Cells.Item(467, 1))
$ws.Cells.Item(467, 1).Value = 44541
$ws.Cells.Item(467, 2).Value = 2
$ws.Cells.Item(467, 3).Value = 25
$ws.Cells.Item(467, 4).Value = 721.917412647993

$dateStyleSource.Copy($ws.Cells.Item(468, 1))
$ws.Cells.Item(468, 1).Value = 44542
$ws.Cells.Item(468, 2).Value = 1
$ws.Cells.Item(468, 3).Value = 17
$ws.Cells.Item(468, 4).Value = 490.9038406006353

$dateStyleSource.Copy($ws.Cells.Item(469, 1))
$ws.Cells.Item(469, 1).Value = 44543
$ws.Cells.Item(469, 2).Value = 1
$ws.Cells.Item(469, 3).Value = 16
$ws.Cells.Item(469, 4).Value = 462.0271440947156

$dateStyleSource.Copy($ws.Cells.Item(470, 1))
$ws.Cells.Item(470, 1).Value = 44544
$ws.Cells.Item(470, 2).Value = 0
$ws.Cells.Item(470, 3).Value = 14
$ws.Cells.Item(470, 4).Value = 404.2737510828762

$dateStyleSource.Copy($ws.Cells.Item(471, 1))
$ws.Cells.Item(471, 1).Value = 44545
$ws.Cells.Item(471, 2).Value = 0
$ws.Cells.Item(471, 3).Value = 13
$ws.Cells.Item(471, 4).Value = 375.3970545769564

$dateStyleSource.Copy($ws.Cells.Item(472, 1))
$ws.Cells.Item(472, 1).Value = 44546
$ws.Cells.Item(472, 2).Value = 0
$ws.Cells.Item(472, 3).Value = 9
$ws.Cells.Item(472, 4).Value = 259.8902685532775

$dateStyleSource.Copy($ws.Cells.Item(473, 1))
$ws.Cells.Item(473, 1).Value = 44547
$ws.Cells.Item(473, 2).Value = 1
$ws.Cells.Item(473, 3).Value = 5
$ws.Cells.Item(473, 4).Value = 144.3834825295986

$dateStyleSource.Copy($ws.Cells.Item(474, 1))
$ws.Cells.Item(474, 1).Value = 44548
$ws.Cells.Item(474, 2).Value = 0
$ws.Cells.Item(474, 3).Value = 3
$ws.Cells.Item(474, 4).Value = 86.63008951775916

$dateStyleSource.Copy($ws.Cells.Item(475, 1))
$ws.Cells.Item(475, 1).Value = 44550
$ws.Cells.Item(475, 2).Value = 0
$ws.Cells.Item(475, 3).Value = 2
$ws.Cells.Item(475, 4).Value = 57.75339301183945

$dateStyleSource.Copy($ws.Cells.Item(476, 1))
$ws.Cells.Item(476, 1).Value = 44551
$ws.Cells.Item(476, 2).Value = 0
$ws.Cells.Item(476, 3).Value = 1
$ws.Cells.Item(476, 4).Value = 28.87669650591972

$dateStyleSource.Copy($ws.Cells.Item(477, 1))
$ws.Cells.Item(477, 1).Value = 44552
$ws.Cells.Item(477, 2).Value = 0
$ws.Cells.Item(477, 3).Value = 1
$ws.Cells.Item(477, 4).Value = 28.87669650591972

$dateStyleSource.Copy($ws.Cells.Item(478, 1))
$ws.Cells.Item(478, 1).Value = 44553
$ws.Cells.Item(478, 2).Value = 2
$ws.Cells.Item(478, 3).Value = 3
$ws.Cells.Item(478, 4).Value = 86.63008951775916

$dateStyleSource.Copy($ws.Cells.Item(479, 1))
$ws.Cells.Item(479, 1).Value = 44554
$ws.Cells.Item(479, 2).Value = 1
$ws.Cells.Item(479, 3).Value = 4
$ws.Cells.Item(479, 4).Value = 115.5067860236789

$dateStyleSource.Copy($ws.Cells.Item(480, 1))
$ws.Cells.Item(480, 1).Value = 44555
$ws.Cells.Item(480, 2).Value = 2
$ws.Cells.Item(480, 3).Value = 5
$ws.Cells.Item(480, 4).Value = 144.3834825295986

$dateStyleSource.Copy($ws.Cells.Item(481, 1))
$ws.Cells.Item(481, 1).Value = 44556
$ws.Cells.Item(481, 2).Value = 2
$ws.Cells.Item(481, 3).Value = 7
$ws.Cells.Item(481, 4).Value = 202.1368755414381

$dateStyleSource.Copy($ws.Cells.Item(482, 1))
$ws.Cells.Item(482, 1).Value = 44557
$ws.Cells.Item(482, 2).Value = 1
$ws.Cells.Item(482, 3).Value = 8
$ws.Cells.Item(482, 4).Value = 231.0135720473578

$dateStyleSource.Copy($ws.Cells.Item(483, 1))
$ws.Cells.Item(483, 1).Value = 44558
$ws.Cells.Item(483, 2).Value = 5
$ws.Cells.Item(483, 3).Value = 13
$ws.Cells.Item(483, 4).Value = 375.3970545769564

$dateStyleSource.Copy($ws.Cells.Item(484, 1))
$ws.Cells.Item(484, 1).Value = 44559
$ws.Cells.Item(484, 2).Value = 3
$ws.Cells.Item(484, 3).Value = 16
$ws.Cells.Item(484, 4).Value = 462.0271440947156

$dateStyleSource.Copy($ws.Cells.Item(485, 1))
$ws.Cells.Item(485, 1).Value = 44560
$ws.Cells.Item(485, 2).Value = 2
$ws.Cells.Item(485, 3).Value = 16
$ws.Cells.Item(485, 4).Value = 462.0271440947156

$dateStyleSource.Copy($ws.Cells.Item(486, 1))
$ws.Cells.Item(486, 1).Value = 44561
$ws.Cells.Item(486, 2).Value = 3
$ws.Cells.Item(486, 3).Value = 18
$ws.Cells.Item(486, 4).Value = 519.780537106555

$dateStyleSource.Copy($ws.Cells.Item(487, 1))
$ws.Cells.Item(487, 1).Value = 44562
$ws.Cells.Item(487, 2).Value = 2
$ws.Cells.Item(487, 3).Value = 18
$ws.Cells.Item(487, 4).Value = 519.780537106555

$dateStyleSource.Copy($ws.Cells.Item(488, 1))
$ws.Cells.Item(488, 1).Value = 44563
$ws.Cells.Item(488, 2).Value = 7
$ws.Cells.Item(488, 3).Value = 23
$ws.Cells.Item(488, 4).Value = 664.1640196361536

$dateStyleSource.Copy($ws.Cells.Item(489, 1))
$ws.Cells.Item(489, 1).Value = 44564
$ws.Cells.Item(489, 2).Value = 9
$ws.Cells.Item(489, 3).Value = 31
$ws.Cells.Item(489, 4).Value = 895.1775916835113

$dateStyleSource.Copy($ws.Cells.Item(490, 1))
$ws.Cells.Item(490, 1).Value = 44565
$ws.Cells.Item(490, 2).Value = 2
$ws.Cells.Item(490, 3).Value = 28
$ws.Cells.Item(490, 4).Value = 808.5475021657523

$dateStyleSource.Copy($ws.Cells.Item(491, 1))
$ws.Cells.Item(491, 1).Value = 44566
$ws.Cells.Item(491, 2).Value = 5
$ws.Cells.Item(491, 3).Value = 30
$ws.Cells.Item(491, 4).Value = 866.3008951775917
